$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.899.68"
$ws.Range("E2").Value = "  -0.58%  "

$ws.Range("D3").Value = "2.356.93"
$ws.Range("E3").Value = "  -0.38%  "

$ws.Range("E4").Value = "  +0.29%  "

$ws.Range("D5").Value = "543.36"
$ws.Range("E5").Value = "  -0.14%  "

$ws.Range("D6").Value = "134.64"
$ws.Range("E6").Value = "  -0.59%  "

$ws.Range("E7").Value = "  +0.16%  "

$ws.Range("D8").Value = "0.562"
$ws.Range("E8").Value = "  +4.89%  "

$ws.Range("E9").Value = "  -0.38%  "

$ws.Range("E10").Value = "  +2.31%  "

$ws.Range("E11").Value = "  -1.36%  "

$ws.Range("D12").Value = "0.355"
$ws.Range("E12").Value = "  +0.01%  "

$ws.Range("D13").Value = "2.776.19"
$ws.Range("E13").Value = "  -0.36%  "

$ws.Range("D14").Value = "23.72"
$ws.Range("E14").Value = "  +0.05%  "

$ws.Range("D15").Value = "57.846.10"
$ws.Range("E15").Value = "  -0.50%  "

$ws.Range("E16").Value = "  +0.39%  "

$ws.Range("D17").Value = "2.356.04"
$ws.Range("E17").Value = "  -0.77%  "

$ws.Range("D18").Value = "10.75"
$ws.Range("E18").Value = "  +2.39%  "

$ws.Range("D19").Value = "331.01"
$ws.Range("E19").Value = "  -2.28%  "

$ws.Range("E20").Value = "  +1.54%  "

$ws.Range("D21").Value = "6.71"
$ws.Range("E21").Value = "  -3.02%  "

$ws.Range("E22").Value = "  +0.15%  "

$ws.Range("D23").Value = "62.58"
$ws.Range("E23").Value = "  +1.00%  "

$ws.Range("D24").Value = "0.166"
$ws.Range("E24").Value = "  -1.95%  "

$ws.Range("E25").Value = "  +0.38%  "

$ws.Range("E26").Value = "  -0.92%  "

$ws.Range("D27").Value = "1.35"
$ws.Range("E27").Value = "  -3.46%  "

$ws.Range("E28").Value = "  -0.15%  "

$ws.Range("D29").Value = "170.18"
$ws.Range("E29").Value = "  -2.31%  "

$ws.Range("D30").Value = "0.0₃0739"
$ws.Range("E30").Value = "  +0.18%  "

$ws.Range("D31").Value = "6.14"
$ws.Range("E31").Value = "  +0.00%  "

$ws.Range("E32").Value = "  +2.60%  "

$ws.Range("D33").Value = "18.41"
$ws.Range("E33").Value = "  -0.85%  "

$ws.Range("E34").Value = "  +0.03%  "

$ws.Range("D35").Value = "4.22"
$ws.Range("E35").Value = "  +3.09%  "

$ws.Range("E36").Value = "  +0.03%  "

$ws.Range("E37").Value = "  -2.18%  "

$ws.Range("E38").Value = "  +0.26%  "

$ws.Range("D39").Value = "39.38"
$ws.Range("E39").Value = "  -0.11%  "

$ws.Range("D40").Value = "142.58"
$ws.Range("E40").Value = "  -5.28%  "

$ws.Range("E41").Value = "  +0.16%  "

$ws.Range("E42").Value = "  +0.83%  "

$ws.Range("D43").Value = "288.61"
$ws.Range("E43").Value = "  +1.32%  "

$ws.Range("D44").Value = "0.0944"
$ws.Range("E44").Value = "  +1.67%  "

$ws.Range("E45").Value = "  +0.62%  "

$ws.Range("D46").Value = "19.14"
$ws.Range("E46").Value = "  -0.20%  "

$ws.Range("E47").Value = "  +0.75%  "

$ws.Range("E48").Value = "  +2.09%  "

$ws.Range("D49").Value = "0.385"
$ws.Range("E49").Value = "  +0.51%  "

$ws.Range("D50").Value = "17.47"
$ws.Range("E50").Value = "  -0.50%  "

$ws.Range("D51").Value = "11.07"
$ws.Range("E51").Value = "  +1.67%  "
